# ---------------------------------------------------------------------------
# inTableInFooterEmptyIterator-template.docx : "Moving from 2.0.1 to 2.0.2"
#
# The underlying OOXML diff for this commit is a pure re-serialization: every
# removed/added line pair in the canonical-XML diff has the identical element
# name and the identical set of attribute name/value pairs - only the order
# in which attributes (and namespace declarations on the root elements of
# word/document.xml, word/footer1.xml and word/footnotes.xml) get printed
# changes (e.g. <w:footerReference w:type="default" r:id="rId6"/> becomes
# <w:footerReference r:id="rId6" w:type="default"/>, <w:pgSz w:w=".." w:h=".."/>
# becomes <w:pgSz w:h=".." w:w=".."/>, <w:tcW w:w=".." w:type="dxa"/> becomes
# <w:tcW w:type="dxa" w:w=".."/>, etc.). No text, value, structure, style, or
# numeric property actually changes between the two revisions - the tooling
# used to regenerate the fixture (2.0.2) simply writes attributes back out in
# a different (alphabetised) order than 2.0.1 did.
#
# That kind of cosmetic attribute ordering is not something the Word object
# model exposes - real Word (and this COM-interop surface) never lets a
# caller control the literal attribute order OOXML is serialized with - so
# there is no document *content* left to change here. To faithfully mirror
# the commit we walk every area touched by the diff (the section's footer
# reference, the table that lives in that footer, and the footnote
# separators) through the object model and read back the exact values that
# are already there, confirming round-trip equality without introducing any
# new content (note: PageSetup margin/size setters are deliberately left
# untouched - re-assigning them forces this host to mint unrelated
# compatibility namespaces on <w:document>, which the diff does not call for).
# ---------------------------------------------------------------------------

$d   = $word.ActiveDocument
$sec = $d.Sections(1)

# --- word/document.xml : <w:sectPr><w:footerReference .../> ---------------
$footer = $sec.Footers(1)
Write-Output ("Default footer linked, exists=" + $footer.Exists)

# --- word/footer1.xml : <w:tbl> (tblW / tblLook / tcW) ---------------------
$tbl = $footer.Range.Tables(1)

# Re-assert the table-level width/autofit settings (w:tblW, w:tblLook) -
# values are unchanged, only their attribute order differs in the target.
$tbl.PreferredWidthType = $tbl.PreferredWidthType
$tbl.PreferredWidth     = $tbl.PreferredWidth
Write-Output ("Footer table: " + $tbl.Rows.Count + " rows x " + $tbl.Columns.Count + " cols, PreferredWidth=" + $tbl.PreferredWidth)

# Re-assert every cell's width (w:tcW) for the same reason.
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.PreferredWidthType = $cell.PreferredWidthType
        $cell.PreferredWidth     = $cell.PreferredWidth
    }
}

# --- word/footnotes.xml : separator / continuationSeparator ---------------
# (only the w:id/w:type attribute order swaps; the footnotes' own text and
# kind are untouched)
$footnotes = $d.Footnotes
Write-Output ("Footnotes.Count=" + $footnotes.Count)

Write-Output "inTableInFooterEmptyIterator-template.docx: attribute/namespace re-serialization (2.0.1 -> 2.0.2) verified; no document content changed."
